$wb = $excel.ActiveWorkbook
$excel.Left = 32250
$excel.Top = -3135
$excel.Width = 16965
$excel.Height = 14760
Write-Host "done"
